$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 8267
$ws.Range("E2").Value = 451
$ws.Range("F2").Value = 459
$ws.Range("G2").Value = 531
$ws.Range("H2").Value = 409
$ws.Range("I2").Value = 399
$ws.Range("J2").Value = 10
$ws.Range("K2").Value = 5121
$ws.Range("L2").Value = 1296
$ws.Range("M2").Value = 3825
$ws.Range("N2").Value = 3764
$ws.Range("O2").Value = 61
$ws.Range("P2").Value = 242
$ws.Range("Q2").Value = 371
$ws.Range("R2").Value = -101
$ws.Range("S2").Value = -35
$ws.Range("T2").Value = 205
$ws.Range("U2").Value = 165
$ws.Range("V2").Value = 522
$ws.Range("W2").Value = 5.46
$ws.Range("X2").Value = 4.95
$ws.Range("Y2").Value = 11.07
$ws.Range("Z2").Value = 8.32
$ws.Range("AA2").Value = 33.87
$ws.Range("AB2").Value = 1413.11
$ws.Range("AC2").Value = 8243
$ws.Range("AD2").Value = 7.12
$ws.Range("AE2").Value = 77767
$ws.Range("AF2").Value = 0.75
$ws.Range("AG2").Value = 1750
$ws.Range("AH2").Value = 2.98
$ws.Range("AI2").Value = 21.23
$ws.Range("AJ2").Value = 4840000
$ws.Range("D3").Value = 7439
$ws.Range("E3").Value = 488
$ws.Range("F3").Value = 488
$ws.Range("G3").Value = 476
$ws.Range("H3").Value = 356
$ws.Range("I3").Value = 337
$ws.Range("J3").Value = 19
$ws.Range("K3").Value = 5076
$ws.Range("L3").Value = 1010
$ws.Range("M3").Value = 4066
$ws.Range("N3").Value = 3985
$ws.Range("O3").Value = 81
$ws.Range("P3").Value = 242
$ws.Range("Q3").Value = 424
$ws.Range("R3").Value = -414
$ws.Range("S3").Value = -199
$ws.Range("T3").Value = 238
$ws.Range("U3").Value = 186
$ws.Range("V3").Value = 412
$ws.Range("W3").Value = 6.57
$ws.Range("X3").Value = 4.78
$ws.Range("Y3").Value = 8.699999999999999
$ws.Range("Z3").Value = 6.98
$ws.Range("AA3").Value = 24.84
$ws.Range("AB3").Value = 1519.87
$ws.Range("AC3").Value = 6968
$ws.Range("AD3").Value = 6.92
$ws.Range("AE3").Value = 82337
$ws.Range("AF3").Value = 0.59
$ws.Range("AG3").Value = 1750
$ws.Range("AH3").Value = 3.63
$ws.Range("AI3").Value = 25.12
$ws.Range("AJ3").Value = 4840000
$ws.Range("D4").Value = 6474
$ws.Range("E4").Value = 443
$ws.Range("F4").Value = 443
$ws.Range("G4").Value = 448
$ws.Range("H4").Value = 341
$ws.Range("I4").Value = 328
$ws.Range("J4").Value = 12
$ws.Range("K4").Value = 5320
$ws.Range("L4").Value = 1020
$ws.Range("M4").Value = 4300
$ws.Range("N4").Value = 4216
$ws.Range("O4").Value = 84
$ws.Range("P4").Value = 242
$ws.Range("Q4").Value = 638
$ws.Range("R4").Value = -123
$ws.Range("S4").Value = -196
$ws.Range("T4").Value = 201
$ws.Range("U4").Value = 437
$ws.Range("V4").Value = 316
$ws.Range("W4").Value = 6.84
$ws.Range("X4").Value = 5.26
$ws.Range("Y4").Value = 8
$ws.Range("Z4").Value = 6.55
$ws.Range("AA4").Value = 23.71
$ws.Range("AB4").Value = 1621.66
$ws.Range("AC4").Value = 6782
$ws.Range("AD4").Value = 8.76
$ws.Range("AE4").Value = 87108
$ws.Range("AF4").Value = 0.68
$ws.Range("AG4").Value = 2000
$ws.Range("AH4").Value = 3.37
$ws.Range("AI4").Value = 29.49
$ws.Range("AJ4").Value = 4840000
$ws.Range("D5").Value = 7037
$ws.Range("E5").Value = 300
$ws.Range("F5").Value = 300
$ws.Range("G5").Value = 354
$ws.Range("H5").Value = 271
$ws.Range("I5").Value = 265
$ws.Range("J5").Value = 6
$ws.Range("K5").Value = 5481
$ws.Range("L5").Value = 1053
$ws.Range("M5").Value = 4428
$ws.Range("N5").Value = 4346
$ws.Range("O5").Value = 82
$ws.Range("P5").Value = 242
$ws.Range("Q5").Value = 393
$ws.Range("R5").Value = -412
$ws.Range("S5").Value = -218
$ws.Range("T5").Value = 235
$ws.Range("U5").Value = 158
$ws.Range("V5").Value = 242
$ws.Range("W5").Value = 4.26
$ws.Range("X5").Value = 3.85
$ws.Range("Y5").Value = 6.19
$ws.Range("Z5").Value = 5.02
$ws.Range("AA5").Value = 23.78
$ws.Range("AB5").Value = 1691.25
$ws.Range("AC5").Value = 5475
$ws.Range("AD5").Value = 12.77
$ws.Range("AE5").Value = 91390
$ws.Range("AF5").Value = 0.76
$ws.Range("AG5").Value = 2000
$ws.Range("AH5").Value = 2.86
$ws.Range("AI5").Value = 35.89
$ws.Range("AJ5").Value = 4840000
$ws.Range("D6").Value = 7667
$ws.Range("E6").Value = 238
$ws.Range("F6").Value = 238
$ws.Range("G6").Value = 157
$ws.Range("H6").Value = 120
$ws.Range("I6").Value = 111
$ws.Range("K6").Value = 5353
$ws.Range("L6").Value = 916
$ws.Range("M6").Value = 4437
$ws.Range("N6").Value = 4350
$ws.Range("P6").Value = 242
$ws.Range("Q6").Value = 214
$ws.Range("R6").Value = -395
$ws.Range("S6").Value = -103
$ws.Range("T6").Value = 329
$ws.Range("U6").Value = -115
$ws.Range("V6").Value = 239
$ws.Range("W6").Value = 3.11
$ws.Range("X6").Value = 1.57
$ws.Range("Y6").Value = 2.54
$ws.Range("Z6").Value = 2.22
$ws.Range("AA6").Value = 20.65
$ws.Range("AB6").Value = 1711.56
$ws.Range("AC6").Value = 2284
$ws.Range("AD6").Value = 22.94
$ws.Range("AE6").Value = 91476
$ws.Range("AF6").Value = 0.57
$ws.Range("AG6").Value = 2000
$ws.Range("AH6").Value = 3.82
$ws.Range("AI6").Value = 86.02
$ws.Range("AJ6").Value = 4840000
$ws.Range("D7").Value = 7560
$ws.Range("E7").Value = 467
$ws.Range("G7").Value = 536
$ws.Range("H7").Value = 394
$ws.Range("I7").Value = 383
$ws.Range("K7").Value = 5826
$ws.Range("L7").Value = 1087
$ws.Range("M7").Value = 4738
$ws.Range("N7").Value = 4674
$ws.Range("P7").Value = 242
$ws.Range("Q7").Value = 1128
$ws.Range("R7").Value = -612
$ws.Range("S7").Value = -176
$ws.Range("T7").Value = 230
$ws.Range("U7").Value = 893
$ws.Range("W7").Value = 6.18
$ws.Range("X7").Value = 5.22
$ws.Range("Y7").Value = 8.49
$ws.Range("Z7").Value = 7.06
$ws.Range("AA7").Value = 22.94
$ws.Range("AC7").Value = 7913
$ws.Range("AD7").Value = 6.8
$ws.Range("AE7").Value = 98292
$ws.Range("AF7").Value = 0.55
$ws.Range("AG7").Value = 2125
$ws.Range("AH7").Value = 3.95
$ws.Range("AI7").Value = 26.85
$ws.Range("D8").Value = 7870
$ws.Range("E8").Value = 529
$ws.Range("G8").Value = 577
$ws.Range("H8").Value = 438
$ws.Range("I8").Value = 428
$ws.Range("K8").Value = 6160
$ws.Range("L8").Value = 1096
$ws.Range("M8").Value = 5064
$ws.Range("N8").Value = 5036
$ws.Range("P8").Value = 242
$ws.Range("Q8").Value = 532
$ws.Range("R8").Value = -314
$ws.Range("S8").Value = -118
$ws.Range("T8").Value = 188
$ws.Range("U8").Value = 372
$ws.Range("W8").Value = 6.72
$ws.Range("X8").Value = 5.57
$ws.Range("Y8").Value = 8.800000000000001
$ws.Range("Z8").Value = 7.32
$ws.Range("AA8").Value = 21.64
$ws.Range("AC8").Value = 8833
$ws.Range("AD8").Value = 6.09
$ws.Range("AE8").Value = 105905
$ws.Range("AF8").Value = 0.51
$ws.Range("AG8").Value = 2375
$ws.Range("AH8").Value = 4.41
$ws.Range("AI8").Value = 26.89
$ws.Range("D9").Value = 7662
$ws.Range("E9").Value = 504
$ws.Range("G9").Value = 550
$ws.Range("H9").Value = 418
$ws.Range("I9").Value = 407
$ws.Range("K9").Value = 6342
$ws.Range("L9").Value = 1100
$ws.Range("M9").Value = 5242
$ws.Range("N9").Value = 5126
$ws.Range("P9").Value = 242
$ws.Range("Q9").Value = 543
$ws.Range("R9").Value = -314
$ws.Range("S9").Value = -126
$ws.Range("T9").Value = 188
$ws.Range("U9").Value = 440
$ws.Range("W9").Value = 6.58
$ws.Range("X9").Value = 5.46
$ws.Range("Y9").Value = 8.01
$ws.Range("Z9").Value = 6.69
$ws.Range("AA9").Value = 20.97
$ws.Range("AC9").Value = 8409
$ws.Range("AD9").Value = 6.4
$ws.Range("AE9").Value = 107797
$ws.Range("AF9").Value = 0.5
$ws.Range("AG9").Value = 2500
$ws.Range("AH9").Value = 4.65
$ws.Range("AI9").Value = 29.73
